# "mejoras para subir archivos"
# Adds a new "Organizacion" column (Z) with header + value "ParBros",
# and changes the "Posicion Politica" value in F2 from the numeric 3
# to the text "centro". Also updates the saved cursor/selection position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Give the new header cell (Z1) the same formatting as the other header
# cells (copy style from Y1, the last existing header), then set values.
$ws.Range("Y1").Copy()
$ws.Range("Z1").PasteSpecial(-4122)

$ws.Range("Z1").Value = "Organizacion"
$ws.Range("Z2").Value = "ParBros"

$ws.Range("F2").Value = "centro"

# Move the selection / active cell to F3 (also clears the old frozen
# scroll position at F1).
$ws.Range("F3").Select()
